# Update the "取得日時" (retrieved datetime) timestamp for all existing
# rows on the "ランサーズ" sheet from 2026-01-27 12:43:44 to 2026-01-27 12:59:18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2026-01-27 12:43:44"
$newValue = "2026-01-27 12:59:18"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
